$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.124.96'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.297.71'
$ws.Range("E3").Value = '  -5.34%  '
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '565.40'
$ws.Range("E5").Value = '  -3.63%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '127.89'
$ws.Range("E6").Value = '  -3.34%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.297.75'
$ws.Range("E8").Value = '  -5.32%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.479'
$ws.Range("E9").Value = '  -0.74%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.32'
$ws.Range("E10").Value = '  -4.50%  '
$ws.Range("E11").Value = '  -4.26%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.374'
$ws.Range("E12").Value = '  -3.25%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.863.11'
$ws.Range("E13").Value = '  -5.38%  '
$ws.Range("E14").Value = '  -0.47%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.301.12'
$ws.Range("E15").Value = '  -5.32%  '
$ws.Range("E16").Value = '  -5.45%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.189.24'
$ws.Range("E17").Value = '  -4.65%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '24.38'
$ws.Range("E18").Value = '  -0.35%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.65'
$ws.Range("E19").Value = '  -1.58%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.32'
$ws.Range("E20").Value = '  -1.29%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '8.95'
$ws.Range("E21").Value = '  -10.58%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '354.30'
$ws.Range("E22").Value = '  -7.89%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.552'
$ws.Range("E23").Value = '  -4.24%  '
$ws.Range("E24").Value = '  -0.05%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.430.31'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '69.12'
$ws.Range("E26").Value = '  -7.04%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.999'
$ws.Range("E28").Value = '  -0.03%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.17'
$ws.Range("E29").Value = '  -0.27%  '
$ws.Range("E30").Value = '  -1.86%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.80'
$ws.Range("E31").Value = '  -1.89%  '
$ws.Range("E32").Value = '  -6.03%  '
$ws.Range("E33").Value = '  -0.04%  '
$ws.Range("E34").Value = '  -2.68%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.326.15'
$ws.Range("E35").Value = '  -5.33%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '22.50'
$ws.Range("E36").Value = '  -2.23%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.24'
$ws.Range("E37").Value = '  -0.69%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.76'
$ws.Range("E38").Value = '  -1.03%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '162.70'
$ws.Range("E39").Value = '  -0.37%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.47'
$ws.Range("E40").Value = '  -3.43%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0755'
$ws.Range("E41").Value = '  -3.00%  '
$ws.Range("E42").Value = '  -0.09%  '
$ws.Range("E43").Value = '  +0.29%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.05'
$ws.Range("E44").Value = '  -1.28%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.741'
$ws.Range("E45").Value = '  -7.65%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.11'
$ws.Range("E46").Value = '  -4.47%  '
$ws.Range("E47").Value = '  -4.88%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '22.17'
$ws.Range("E48").Value = '  -7.88%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.67'
$ws.Range("E49").Value = '  -1.11%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.851'
$ws.Range("E50").Value = '  -9.01%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '21.13'
$ws.Range("E51").Value = '  +2.37%  '
